$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.922.19'
$ws.Range("E2").Value = '  +8.17%  '

$ws.Range("D3").Value = '1.768.59'
$ws.Range("E3").Value = '  +6.40%  '

$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '''317.09'
$ws.Range("E5").Value = '  +2.84%  '

$ws.Range("D6").Value = '''0.9974'
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("D7").Value = '''0.3835'
$ws.Range("E7").Value = '  +3.41%  '

$ws.Range("D8").Value = '''0.3618'
$ws.Range("E8").Value = '  +5.39%  '

$ws.Range("D9").Value = '''50.68'
$ws.Range("E9").Value = '  +6.17%  '

$ws.Range("D10").Value = '''1.235'
$ws.Range("E10").Value = '  +5.93%  '

$ws.Range("D11").Value = '''0.07707'
$ws.Range("E11").Value = '  +6.71%  '

$ws.Range("D12").Value = '''0.9979'
$ws.Range("E12").Value = '  -0.10%  '

$ws.Range("D13").Value = '''21.77'
$ws.Range("E13").Value = '  +6.25%  '

$ws.Range("D14").Value = '''6.502'
$ws.Range("E14").Value = '  +8.15%  '

$ws.Range("D15").Value = '''7.118'
$ws.Range("E15").Value = '  +5.92%  '

$ws.Range("D16").Value = '1.766.98'
$ws.Range("E16").Value = '  +6.19%  '

$ws.Range("D17").Value = '''0.00001160'
$ws.Range("E17").Value = '  +6.09%  '

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '''0.06791'
$ws.Range("E18").Value = '  +1.23%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '''0.9974'
$ws.Range("E19").Value = '  +0.13%  '

$ws.Range("D20").Value = '''87.29'
$ws.Range("E20").Value = '  +7.19%  '

$ws.Range("D21").Value = '''17.77'
$ws.Range("E21").Value = '  +8.46%  '

$ws.Range("D22").Value = '''6.528'
$ws.Range("E22").Value = '  +6.85%  '

$ws.Range("D23").Value = '''12.83'
$ws.Range("E23").Value = '  +7.36%  '

$ws.Range("D24").Value = '25.862.72'
$ws.Range("E24").Value = '  +8.12%  '

$ws.Range("D25").Value = '''2.420'
$ws.Range("E25").Value = '  +1.35%  '

$ws.Range("D26").Value = '''2.943'
$ws.Range("E26").Value = '  +10.39%  '

$ws.Range("D27").Value = '''20.77'
$ws.Range("E27").Value = '  +6.65%  '

$ws.Range("D28").Value = '''156.56'
$ws.Range("E28").Value = '  +3.26%  '

$ws.Range("D29").Value = '1.963.19'
$ws.Range("E29").Value = '  +6.30%  '

$ws.Range("D30").Value = '''134.43'
$ws.Range("E30").Value = '  +5.99%  '

$ws.Range("D31").Value = '''1.232'
$ws.Range("E31").Value = '  +26.39%  '

$ws.Range("D32").Value = '''7.260'
$ws.Range("E32").Value = '  +14.64%  '

$ws.Range("D33").Value = '''4.239'
$ws.Range("E33").Value = '  +3.25%  '

$ws.Range("D34").Value = '''14.23'
$ws.Range("E34").Value = '  +16.24%  '

$ws.Range("D35").Value = '''1.822'
$ws.Range("E35").Value = '  +4.68%  '

$ws.Range("D36").Value = '''0.08794'
$ws.Range("E36").Value = '  +5.01%  '

$ws.Range("D37").Value = '''5.724'
$ws.Range("E37").Value = '  +8.11%  '

$ws.Range("D38").Value = '''0.06787'
$ws.Range("E38").Value = '  +7.14%  '

$ws.Range("D39").Value = '''0.02508'
$ws.Range("E39").Value = '  +8.56%  '

$ws.Range("D40").Value = '''9.418'
$ws.Range("E40").Value = '  +5.45%  '

$ws.Range("D41").Value = '''0.2262'
$ws.Range("E41").Value = '  +9.25%  '

$ws.Range("D42").Value = '''1.299'
$ws.Range("E42").Value = '  +1.47%  '

$ws.Range("D43").Value = '''0.6604'
$ws.Range("E43").Value = '  +8.66%  '

$ws.Range("D44").Value = '''14.47'
$ws.Range("E44").Value = '  +9.42%  '

$ws.Range("D45").Value = '''0.9969'
$ws.Range("E45").Value = '  +0.15%  '

$ws.Range("D46").Value = '''0.6376'
$ws.Range("E46").Value = '  +7.67%  '

$ws.Range("D47").Value = '''3.914'
$ws.Range("E47").Value = '  +2.45%  '

$ws.Range("D48").Value = '''2.176'
$ws.Range("E48").Value = '  +9.27%  '

$ws.Range("D49").Value = '''133.42'
$ws.Range("E49").Value = '  +5.19%  '

$ws.Range("D50").Value = '''0.07509'
$ws.Range("E50").Value = '  +6.20%  '

$ws.Range("D51").Value = '''81.15'
$ws.Range("E51").Value = '  +7.17%  '
